# "change code vehicle, teacher, combobox and fix bug"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (vehicle #1) ---
$ws.Range("B6").Value = "77B0999   "
$ws.Range("F6").Value = "Toyota         "

# G6 needs to become the text "2023      " (not the number 2023). A plain
# assignment gets auto-coerced to a number by Excel's type inference, so
# copy the already-text "2023      " from G7 (values only) to keep G6's
# cell type as text/shared-string while preserving G6's own style.
$ws.Range("G7").Copy()
$ws.Range("G6").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false

$ws.Range("H6").Value = "Toyota         "
$ws.Range("I6").Value = "B1 "
$ws.Range("M6").Value = "Quang Ð?t                     "

# --- Row 7 (vehicle #2) ---
$ws.Range("E7").Value = "RR             "
$ws.Range("F7").Value = "RR             "
$ws.Range("H7").Value = "RR             "
$ws.Range("I7").Value = "B1 "
$ws.Range("M7").Value = "Quoc Anh                      "

# --- Updated submission timestamp ---
$ws.Range("K8").Value = 45183.8613597222

# --- Column M (13) combobox width tweak ---
$ws.Columns.Item(13).ColumnWidth = 10.0246734619141
